$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "58.079.26"
Set-TextValue "E2" "  -0.34%  "
Set-TextValue "D3" "2.587.13"
Set-TextValue "E3" "  -1.85%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "518.57"
Set-TextValue "E5" "  -0.31%  "
Set-TextValue "D6" "143.73"
Set-TextValue "E6" "  +0.49%  "
Set-TextValue "D7" "0.997"
Set-TextValue "E7" "  -0.24%  "
Set-TextValue "D8" "0.567"
Set-TextValue "E8" "  -0.07%  "
Set-TextValue "D9" "2.604.54"
Set-TextValue "E9" "  -1.29%  "
Set-TextValue "D10" "6.69"
Set-TextValue "E10" "  +0.37%  "
Set-TextValue "E11" "  -1.37%  "
Set-TextValue "E13" "  -0.94%  "
Set-TextValue "D14" "3.043.57"
Set-TextValue "E14" "  -1.73%  "
Set-TextValue "D15" "58.030.46"
Set-TextValue "E15" "  -0.41%  "
Set-TextValue "D16" "20.40"
Set-TextValue "E16" "  -1.73%  "
Set-TextValue "E17" "  -1.36%  "
Set-TextValue "D18" "2.606.43"
Set-TextValue "E18" "  -1.15%  "
Set-TextValue "D19" "340.13"
Set-TextValue "E19" "  +1.34%  "
Set-TextValue "D20" "4.32"
Set-TextValue "E20" "  -1.85%  "
Set-TextValue "D21" "10.31"
Set-TextValue "E21" "  -1.17%  "
Set-TextValue "D22" "6.36"
Set-TextValue "E22" "  +1.19%  "
Set-TextValue "E23" "  +0.07%  "
Set-TextValue "D24" "66.09"
Set-TextValue "E24" "  +2.74%  "
Set-TextValue "E25" "  -0.78%  "
Set-TextValue "E26" "  -5.10%  "
Set-TextValue "E27" "  -0.27%  "
Set-TextValue "D28" "2.706.20"
Set-TextValue "E28" "  -1.82%  "
Set-TextValue "D29" "7.01"
Set-TextValue "E29" "  -0.95%  "
Set-TextValue "D30" "0.0₃0748"
Set-TextValue "E30" "  -5.53%  "
Set-TextValue "D31" "0.998"
Set-TextValue "E31" "  -0.07%  "
Set-TextValue "D32" "6.27"
Set-TextValue "E32" "  -4.93%  "
Set-TextValue "E33" "  -0.32%  "
Set-TextValue "D34" "18.78"
Set-TextValue "E34" "  -0.02%  "
Set-TextValue "D35" "149.48"
Set-TextValue "E35" "  -2.12%  "
Set-TextValue "E36" "  -1.75%  "
Set-TextValue "E37" "  -2.49%  "
Set-TextValue "D38" "0.881"
Set-TextValue "E38" "  -2.58%  "
Set-TextValue "D39" "0.840"
Set-TextValue "E39" "  -1.57%  "
Set-TextValue "E40" "  +0.85%  "
Set-TextValue "D41" "35.93"
Set-TextValue "E41" "  -2.35%  "
Set-TextValue "D42" "3.56"
Set-TextValue "E42" "  -1.86%  "
Set-TextValue "D43" "0.996"
Set-TextValue "E43" "  -0.36%  "
Set-TextValue "D44" "272.87"
Set-TextValue "E44" "  +1.14%  "
Set-TextValue "D45" "0.591"
Set-TextValue "E45" "  -1.86%  "
Set-TextValue "D46" "10.65"
Set-TextValue "E46" "  +0.10%  "
Set-TextValue "E47" "  -1.54%  "
Set-TextValue "D48" "18.84"
Set-TextValue "E48" "  -2.34%  "
Set-TextValue "D49" "0.0524"
Set-TextValue "E49" "  -1.92%  "
Set-TextValue "B50" "RenderToken"
Set-TextValue "C50" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D50" "4.67"
Set-TextValue "E50" "  -0.35%  "
Set-TextValue "B51" "Maker"
Set-TextValue "C51" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D51" "1.976.96"
Set-TextValue "E51" "  -2.87%  "
